$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price value is a plain decimal number need their cell
# format set to Text first, otherwise Excel auto-converts the typed value
# into a number (e.g. "0.200" -> 0.2, "1.00" -> 1), same as interactive entry.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "94.995.54"
$ws.Range("E2").Value = "  +1.51%  "
$ws.Range("D3").Value = "3.606.95"
$ws.Range("E3").Value = "  +4.63%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "235.64"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").Value = "655.94"
$ws.Range("E6").Value = "  +4.52%  "
$ws.Range("E7").Value = "  +0.94%  "
$ws.Range("D8").Value = "0.398"
$ws.Range("E8").Value = "  +1.21%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").Value = "0.988"
$ws.Range("E10").Value = "  -1.75%  "
$ws.Range("D11").Value = "3.604.50"
$ws.Range("E11").Value = "  +4.62%  "
$ws.Range("D12").Value = "0.200"
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("D13").Value = "41.94"
$ws.Range("E13").Value = "  -3.73%  "
$ws.Range("D14").Value = "6.33"
$ws.Range("D15").Value = "4.273.88"
$ws.Range("E15").Value = "  +4.17%  "
$ws.Range("D16").Value = "94.903.62"
$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("D18").Value = "3.604.31"
$ws.Range("E18").Value = "  +5.47%  "
$ws.Range("D19").Value = "7.92"
$ws.Range("E19").Value = "  -4.59%  "
$ws.Range("D20").Value = "12.77"
$ws.Range("E20").Value = "  +8.60%  "
$ws.Range("D21").Value = "17.89"
$ws.Range("E21").Value = "  -1.07%  "
$ws.Range("D22").Value = "3.53"
$ws.Range("E22").Value = "  +4.14%  "
$ws.Range("D23").Value = "0.477"
$ws.Range("E23").Value = "  -3.18%  "
$ws.Range("D24").Value = "501.40"
$ws.Range("E24").Value = "  -0.75%  "
$ws.Range("D25").Value = "0.0000195"
$ws.Range("E25").Value = "  +5.66%  "
$ws.Range("D26").Value = "6.57"
$ws.Range("E26").Value = "  -3.39%  "
$ws.Range("D27").Value = "94.56"
$ws.Range("E27").Value = "  -0.81%  "
$ws.Range("D28").Value = "3.799.45"
$ws.Range("E28").Value = "  +4.76%  "
$ws.Range("E29").Value = "  +2.00%  "
$ws.Range("E30").Value = "  +8.81%  "
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("D32").Value = "11.15"
$ws.Range("E32").Value = "  -1.95%  "
$ws.Range("E33").Value = "  -1.28%  "
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  -0.29%  "
$ws.Range("D35").Value = "32.06"
$ws.Range("E35").Value = "  +8.80%  "
$ws.Range("E36").Value = "  -2.48%  "
$ws.Range("D37").Value = "0.555"
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("D38").Value = "561.88"
$ws.Range("E38").Value = "  -1.67%  "
$ws.Range("D39").Value = "7.99"
$ws.Range("E39").Value = "  +5.88%  "
$ws.Range("E40").Value = "  +1.75%  "
$ws.Range("E42").Value = "  +0.39%  "
$ws.Range("D43").Value = "0.908"
$ws.Range("E43").Value = "  -0.71%  "
$ws.Range("D44").Value = "35.35"
$ws.Range("E44").Value = "  +43.59%  "
$ws.Range("D45").Value = "1.72"
$ws.Range("E45").Value = "  +1.11%  "
$ws.Range("D46").Value = "23.68"
$ws.Range("E46").Value = "  -0.17%  "
$ws.Range("E47").Value = "  +0.80%  "
$ws.Range("E48").Value = "  +4.80%  "
$ws.Range("D49").Value = "0.0411"
$ws.Range("E49").Value = "  -3.03%  "
$ws.Range("E50").Value = "  -1.38%  "
$ws.Range("D51").Value = "53.40"
$ws.Range("E51").Value = "  +0.36%  "
